$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# I0 / IF data values for rows 2-31
$values = @{
    2  = @(8, 8)
    3  = @(8, 8)
    4  = @(7, 8)
    5  = @(8, 8)
    6  = @(8, 9)
    7  = @(8, 8)
    8  = @(7, 8)
    9  = @(8, 8)
    10 = @(7, 7)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(8, 8)
    16 = @(8, 8)
    17 = @(8, 8)
    18 = @(7, 8)
    19 = @(7, 7)
    20 = @(7, 7)
    21 = @(8, 8)
    22 = @(8, 8)
    23 = @(8, 8)
    24 = @(8, 8)
    25 = @(9, 9)
    26 = @(8, 8)
    27 = @(8, 8)
    28 = @(5, 5)
    29 = @(5, 5)
    30 = @(5, 5)
    31 = @(3, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Range("I$row").Value = $pair[0]
    $ws.Range("J$row").Value = $pair[1]
}
